$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(252).Insert()

$ws.Cells.Item(252, 1).Value = 7
$ws.Cells.Item(252, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(252, 3).Value = "Ñuble"
$ws.Cells.Item(252, 4).Value = 44932
$ws.Cells.Item(252, 5).Value = 16
$ws.Cells.Item(252, 6).Value = 100112003
$ws.Cells.Item(252, 7).Value = "Ajo"
$ws.Cells.Item(252, 8).Value = "Chino"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 60
$ws.Cells.Item(252, 11).Value = 17000
$ws.Cells.Item(252, 12).Value = 18000
$ws.Cells.Item(252, 13).Value = 17500
$ws.Cells.Item(252, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(252, 15).Value = "China"
$ws.Cells.Item(252, 16).Value = 1750
$ws.Cells.Item(252, 17).Value = 10
$ws.Cells.Item(252, 18).Value = "Hortaliza"
